# Apply the "Add files via upload" edit to FAST_holdings.xlsx:
#   - bump the "as of" date in the confidential disclosure text from
#     2021-05-21 to 2021-05-24
#   - refresh the Weight (D) / Percent Change (E) figures for rows 2-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected, so it must be unprotected before any cell
# contents can change, and re-protected afterwards to preserve the
# original (protected) state of the workbook.
$ws.Unprotect()

# --- Disclosure text (shared string behind cell A13) ---------------------
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) values for rows 2-10 -----------------
$ws.Range("D2").Value = 0.09166615363678114
$ws.Range("E2").Value = 0.008549189740972185

$ws.Range("D3").Value = 0.1060491668183242
$ws.Range("E3").Value = 0.01720586453410888

$ws.Range("D4").Value = 0.119659546284112
$ws.Range("E4").Value = 0.009495409244290931

$ws.Range("D5").Value = 0.1409593609986395
$ws.Range("E5").Value = 0.007549868870698528

$ws.Range("D6").Value = 0.1372956933643795
$ws.Range("E6").Value = 0.003311715192493514

$ws.Range("D7").Value = 0.1469756900603407
$ws.Range("E7").Value = 0.006308257226250058

$ws.Range("D8").Value = 0.1270463529324588
$ws.Range("E8").Value = 0.01491201908738438

$ws.Range("D9").Value = 0.1303480359049642
$ws.Range("E9").Value = 0.008763539175692614

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.009227452443612227

# Restore sheet protection (same settings the workbook shipped with:
# structure/objects/scenarios locked, column & row formatting allowed).
$ws.Protect($null, $true, $true, $true, $false, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true)
